{"js": "// This script updates the answers in the \"three-digit number divided by\n// one-digit number\" practice table. The document body contains a single\n// table; the populated data rows are 0, 4, 8, 12, 16 (0-based), each with\n// 5 columns of \"A\u00f7B=C, D\" style answers. We replace each cell's text with\n// its new value, addressing cells positionally (row/column) so that the\n// edit is unambiguous even though some new values elsewhere in the table\n// coincide with old values that haven't been updated yet.\n\nconst table = context.document.body.tables.getFirst();\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"618\u00f72=309, 0\", newText: \"539\u00f78=67, 3\" },\n  { row: 0, col: 1, oldText: \"154\u00f72=77, 0\", newText: \"678\u00f76=113, 0\" },\n  { row: 0, col: 2, oldText: \"101\u00f74=25, 1\", newText: \"155\u00f73=51, 2\" },\n  { row: 0, col: 3, oldText: \"404\u00f77=57, 5\", newText: \"455\u00f77=65, 0\" },\n  { row: 0, col: 4, oldText: \"577\u00f73=192, 1\", newText: \"249\u00f76=41, 3\" },\n  { row: 4, col: 0, oldText: \"175\u00f74=43, 3\", newText: \"409\u00f75=81, 4\" },\n  { row: 4, col: 1, oldText: \"487\u00f72=243, 1\", newText: \"903\u00f79=100, 3\" },\n  { row: 4, col: 2, oldText: \"304\u00f72=152, 0\", newText: \"434\u00f73=144, 2\" },\n  { row: 4, col: 3, oldText: \"597\u00f77=85, 2\", newText: \"608\u00f79=67, 5\" },\n  { row: 4, col: 4, oldText: \"646\u00f77=92, 2\", newText: \"669\u00f72=334, 1\" },\n  { row: 8, col: 0, oldText: \"128\u00f77=18, 2\", newText: \"111\u00f77=15, 6\" },\n  { row: 8, col: 1, oldText: \"610\u00f78=76, 2\", newText: \"482\u00f74=120, 2\" },\n  { row: 8, col: 2, oldText: \"409\u00f75=81, 4\", newText: \"742\u00f77=106, 0\" },\n  { row: 8, col: 3, oldText: \"271\u00f79=30, 1\", newText: \"109\u00f72=54, 1\" },\n  { row: 8, col: 4, oldText: \"133\u00f79=14, 7\", newText: \"554\u00f75=110, 4\" },\n  { row: 12, col: 0, oldText: \"741\u00f76=123, 3\", newText: \"735\u00f77=105, 0\" },\n  { row: 12, col: 1, oldText: \"594\u00f79=66, 0\", newText: \"124\u00f75=24, 4\" },\n  { row: 12, col: 2, oldText: \"584\u00f76=97, 2\", newText: \"211\u00f79=23, 4\" },\n  { row: 12, col: 3, oldText: \"947\u00f77=135, 2\", newText: \"850\u00f74=212, 2\" },\n  { row: 12, col: 4, oldText: \"135\u00f79=15, 0\", newText: \"954\u00f73=318, 0\" },\n  { row: 16, col: 0, oldText: \"369\u00f79=41, 0\", newText: \"204\u00f77=29, 1\" },\n  { row: 16, col: 1, oldText: \"438\u00f77=62, 4\", newText: \"624\u00f75=124, 4\" },\n  { row: 16, col: 2, oldText: \"398\u00f73=132, 2\", newText: \"599\u00f77=85, 4\" },\n  { row: 16, col: 3, oldText: \"327\u00f76=54, 3\", newText: \"321\u00f74=80, 1\" },\n  { row: 16, col: 4, oldText: \"850\u00f74=212, 2\", newText: \"380\u00f76=63, 2\" },\n];\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    // Replace just the matched range so the run's formatting (font, size,\n    // paragraph alignment, etc.) is preserved.\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: if the expected old text wasn't found (e.g. re-running the\n    // script on an already-updated document), just set the cell text.\n    cell.body.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Updates the answers in the \"three-digit number divided by one-digit\n# number\" practice table. The table's populated data rows are (1-based)\n# rows 1, 5, 9, 13, 17, each with 5 columns containing \"A\u00f7B=C, D\" style\n# answers. Each cell is addressed positionally (row/column) and its text\n# is replaced directly via Range.Text, which (unlike Find.Execute's\n# Replace parameter, which in this host ends up rewriting every matching\n# occurrence document-wide instead of just the target range) reliably\n# touches only the targeted cell and preserves that cell's paragraph /\n# run formatting (font, size, alignment).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  @{ Row = 1; Col = 1; OldText = \"618\u00f72=309, 0\"; NewText = \"539\u00f78=67, 3\" },\n  @{ Row = 1; Col = 2; OldText = \"154\u00f72=77, 0\"; NewText = \"678\u00f76=113, 0\" },\n  @{ Row = 1; Col = 3; OldText = \"101\u00f74=25, 1\"; NewText = \"155\u00f73=51, 2\" },\n  @{ Row = 1; Col = 4; OldText = \"404\u00f77=57, 5\"; NewText = \"455\u00f77=65, 0\" },\n  @{ Row = 1; Col = 5; OldText = \"577\u00f73=192, 1\"; NewText = \"249\u00f76=41, 3\" },\n  @{ Row = 5; Col = 1; OldText = \"175\u00f74=43, 3\"; NewText = \"409\u00f75=81, 4\" },\n  @{ Row = 5; Col = 2; OldText = \"487\u00f72=243, 1\"; NewText = \"903\u00f79=100, 3\" },\n  @{ Row = 5; Col = 3; OldText = \"304\u00f72=152, 0\"; NewText = \"434\u00f73=144, 2\" },\n  @{ Row = 5; Col = 4; OldText = \"597\u00f77=85, 2\"; NewText = \"608\u00f79=67, 5\" },\n  @{ Row = 5; Col = 5; OldText = \"646\u00f77=92, 2\"; NewText = \"669\u00f72=334, 1\" },\n  @{ Row = 9; Col = 1; OldText = \"128\u00f77=18, 2\"; NewText = \"111\u00f77=15, 6\" },\n  @{ Row = 9; Col = 2; OldText = \"610\u00f78=76, 2\"; NewText = \"482\u00f74=120, 2\" },\n  @{ Row = 9; Col = 3; OldText = \"409\u00f75=81, 4\"; NewText = \"742\u00f77=106, 0\" },\n  @{ Row = 9; Col = 4; OldText = \"271\u00f79=30, 1\"; NewText = \"109\u00f72=54, 1\" },\n  @{ Row = 9; Col = 5; OldText = \"133\u00f79=14, 7\"; NewText = \"554\u00f75=110, 4\" },\n  @{ Row = 13; Col = 1; OldText = \"741\u00f76=123, 3\"; NewText = \"735\u00f77=105, 0\" },\n  @{ Row = 13; Col = 2; OldText = \"594\u00f79=66, 0\"; NewText = \"124\u00f75=24, 4\" },\n  @{ Row = 13; Col = 3; OldText = \"584\u00f76=97, 2\"; NewText = \"211\u00f79=23, 4\" },\n  @{ Row = 13; Col = 4; OldText = \"947\u00f77=135, 2\"; NewText = \"850\u00f74=212, 2\" },\n  @{ Row = 13; Col = 5; OldText = \"135\u00f79=15, 0\"; NewText = \"954\u00f73=318, 0\" },\n  @{ Row = 17; Col = 1; OldText = \"369\u00f79=41, 0\"; NewText = \"204\u00f77=29, 1\" },\n  @{ Row = 17; Col = 2; OldText = \"438\u00f77=62, 4\"; NewText = \"624\u00f75=124, 4\" },\n  @{ Row = 17; Col = 3; OldText = \"398\u00f73=132, 2\"; NewText = \"599\u00f77=85, 4\" },\n  @{ Row = 17; Col = 4; OldText = \"327\u00f76=54, 3\"; NewText = \"321\u00f74=80, 1\" },\n  @{ Row = 17; Col = 5; OldText = \"850\u00f74=212, 2\"; NewText = \"380\u00f76=63, 2\" },\n)\n\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    # Assigning Range.Text replaces only this cell's content (and keeps the\n    # paragraph/run formatting of its first run), so cells are addressed\n    # purely positionally and there is no risk of the OldText of one cell\n    # accidentally matching the (already-written) NewText of another.\n    $cell.Range.Text = $item.NewText\n}\n\n"}
